$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.542.33"
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = "'1.993.05"
$ws.Range('E3').Value = '  +6.00%  '
$ws.Range('D4').Value = "'0.9999"
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = "'325.70"
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = "'0.9999"
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').Value = "'0.4678"
$ws.Range('E7').Value = '  +1.57%  '
$ws.Range('D8').Value = "'0.3949"
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('D9').Value = "'46.38"
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('D10').Value = "'0.07931"
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').Value = "'1.001"
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').Value = "'22.90"
$ws.Range('E12').Value = '  +5.19%  '
$ws.Range('D13').Value = "'1.992.85"
$ws.Range('E13').Value = '  +9.57%  '
$ws.Range('D14').Value = "'7.278"
$ws.Range('E14').Value = '  +4.05%  '
$ws.Range('D15').Value = "'5.867"
$ws.Range('E15').Value = '  +3.91%  '
$ws.Range('D16').Value = "'0.07119"
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').Value = "'88.69"
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').Value = "'1.001"
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').Value = "'0.000009981"
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = "'17.41"
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('D21').Value = "'0.9995"
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = "'29.611.71"
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('D23').Value = "'5.533"
$ws.Range('E23').Value = '  +5.70%  '
$ws.Range('D24').Value = "'11.27"
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('D25').Value = "'2.102"
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').Value = "'157.80"
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').Value = "'19.67"
$ws.Range('E27').Value = '  +1.88%  '
$ws.Range('D28').Value = "'5.995"
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = "'120.05"
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('D30').Value = "'1.958"
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').Value = "'0.09451"
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').Value = "'0.9074"
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = "'1.350"
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = "'5.255"
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = "'3.178"
$ws.Range('E35').Value = '  -2.41%  '
$ws.Range('D36').Value = "'0.000003528"
$ws.Range('E36').Value = '  +117.32%  '
$ws.Range('D38').Value = "'1.172"
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').Value = "'0.02117"
$ws.Range('E39').Value = '  +2.19%  '
$ws.Range('D40').Value = "'7.863"
$ws.Range('E40').Value = '  +3.06%  '
$ws.Range('D41').Value = "'0.5753"
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('D42').Value = "'0.1821"
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('D43').Value = "'9.826"
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('D44').Value = "'11.97"
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').Value = "'0.5373"
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = "'2.699"
$ws.Range('E46').Value = '  +6.66%  '
$ws.Range('D47').Value = "'2.164"
$ws.Range('E47').Value = '  -4.67%  '
$ws.Range('D48').Value = "'1.870"
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').Value = "'0.06947"
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('D50').Value = "'114.04"
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('E51').Value = '  +8.32%  '
